$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 626 (shifts existing rows 626.. down to 628..)
$ws.Rows.Item(626).Resize(2).Insert()

# New row 626
$ws.Range("A626").Value = 9
$ws.Range("B626").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C626").Value = "Metropolitana"
$ws.Range("D626").Value = 44491
$ws.Range("E626").Value = 13
$ws.Range("F626").Value = "Fruta"
$ws.Range("G626").Value = 100108
$ws.Range("H626").Value = "Tropicales y subtropicales"
$ws.Range("I626").Value = 100108006
$ws.Range("J626").Value = "Plátano"
$ws.Range("K626").Value = "Sin especificar"
$ws.Range("L626").Value = "Pintón"
$ws.Range("M626").Value = 300
$ws.Range("N626").Value = 24000
$ws.Range("O626").Value = 24000
$ws.Range("P626").Value = 24000
$ws.Range("Q626").Value = "$/caja 20 kilos"
$ws.Range("R626").Value = "Ecuador"
$ws.Range("S626").Value = 1200
$ws.Range("T626").Value = 20

# New row 627
$ws.Range("A627").Value = 9
$ws.Range("B627").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C627").Value = "Metropolitana"
$ws.Range("D627").Value = 44491
$ws.Range("E627").Value = 13
$ws.Range("F627").Value = "Fruta"
$ws.Range("G627").Value = 100108
$ws.Range("H627").Value = "Tropicales y subtropicales"
$ws.Range("I627").Value = 100108006
$ws.Range("J627").Value = "Plátano"
$ws.Range("K627").Value = "Sin especificar"
$ws.Range("L627").Value = "Primera Maduro"
$ws.Range("M627").Value = 380
$ws.Range("N627").Value = 25000
$ws.Range("O627").Value = 25000
$ws.Range("P627").Value = 25000
$ws.Range("Q627").Value = "$/caja 20 kilos"
$ws.Range("R627").Value = "Ecuador"
$ws.Range("S627").Value = 1250
$ws.Range("T627").Value = 20

Write-Host "Done"
